$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Summary")

# These cells hold numeric-looking values that are stored as TEXT
# (shared strings) in the workbook, so assign with a leading apostrophe
# to keep Excel from auto-converting them to numbers.

# Enterprises density (per 1000 people) - SMEs column
$ws.Range("C13").Value = "'2.65"

# Value added to the economy (% of total): Micro, SMEs, MSMEs
$ws.Range("B18").Value = "'34.73"
$ws.Range("C18").Value = "'23.21"
$ws.Range("D18").Value = "'57.94"
